$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 7.7777777
$ws.Range("I11").Value = 7.7777777
$ws.Range("K11").Value = 7.7777777
$ws.Range("M11").Value = 132.2222223

$ws.Range("H51").Value = 9500
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516

$ws.Range("H62").Value = 3723.5

$ws.Range("H65").Value = 3723.5

$ws.Range("H76").Value = 4199.8335
$ws.Range("I76").Value = 4450
$ws.Range("K76").Value = 4450
$ws.Range("M76").Value = -4135

$ws.Range("H79").Value = 4199.8335
$ws.Range("I79").Value = 4450
$ws.Range("K79").Value = 4450
$ws.Range("M79").Value = -3358

$ws.Range("H113").Value = 10475.125
$ws.Range("J113").Value = 10002.5
$ws.Range("L113").Value = 10002.5
$ws.Range("N113").Value = -16510.5

$ws.Range("H118").Value = 1999.5
$ws.Range("J118").Value = 7500
$ws.Range("L118").Value = 22500
$ws.Range("N118").Value = -25814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 66
$ws.Range("N4").ClearContents()

$ws.Range("H45").Value = 4073.8823
$ws.Range("I45").Value = 2250
$ws.Range("J45").Value = 4317.067
$ws.Range("K45").Value = 2250
$ws.Range("L45").Value = 4317.067
$ws.Range("M45").Value = -1873
$ws.Range("N45").Value = -5071.067

$ws.Range("H74").Value = 949.5
$ws.Range("I74").Value = 949.5
$ws.Range("K74").Value = 949.5
$ws.Range("M74").Value = -75.5

$ws.Range("H77").Value = 949.5
$ws.Range("I77").Value = 949.5
$ws.Range("K77").Value = 4747.5
$ws.Range("M77").Value = -379.5

$ws.Range("H88").Value = 3548.375
$ws.Range("I88").Value = 1274.8
$ws.Range("J88").Value = 4581.8184
$ws.Range("K88").Value = 1274.8
$ws.Range("L88").Value = 4581.8184
$ws.Range("M88").Value = -868.8
$ws.Range("N88").Value = -5393.8184

$ws.Range("H91").Value = 3548.375
$ws.Range("I91").Value = 1274.8
$ws.Range("J91").Value = 4581.8184
$ws.Range("K91").Value = 1274.8
$ws.Range("L91").Value = 4581.8184
$ws.Range("M91").Value = 129.2
$ws.Range("N91").Value = -7389.8184

$ws.Range("H97").Value = 1681.4286
$ws.Range("I97").Value = 1100
$ws.Range("K97").Value = 1100
$ws.Range("M97").Value = -604

$ws.Range("H130").Value = 11497
$ws.Range("J130").Value = 11497
$ws.Range("L130").Value = 11497
$ws.Range("N130").Value = -21537

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 5000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -5494

$ws.Range("H22").Value = 1597
$ws.Range("I22").Value = 1597
$ws.Range("K22").Value = 1597
$ws.Range("M22").Value = -1424

$ws.Range("H105").Value = 3638.4285
$ws.Range("I105").Value = 3578.1667
$ws.Range("K105").Value = 3578.1667
$ws.Range("M105").Value = -1831.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1474.75
$ws.Range("I22").Value = 1699.6666
$ws.Range("J22").Value = 1339.8
$ws.Range("K22").Value = 1699.6666
$ws.Range("L22").Value = 1339.8
$ws.Range("M22").Value = -1349.6666
$ws.Range("N22").Value = -2039.8

$ws.Range("H31").Value = 2060.75
$ws.Range("I31").Value = 1434.25
$ws.Range("J31").Value = 2687.25
$ws.Range("K31").Value = 1434.25
$ws.Range("L31").Value = 2687.25
$ws.Range("M31").Value = -1139.25
$ws.Range("N31").Value = -3277.25

$ws.Range("H33").Value = 3886.2
$ws.Range("I33").Value = 3886.2
$ws.Range("K33").Value = 3886.2
$ws.Range("M33").Value = -3507.2

$ws.Range("H34").Value = 2060.75
$ws.Range("I34").Value = 1434.25
$ws.Range("J34").Value = 2687.25
$ws.Range("K34").Value = 1434.25
$ws.Range("L34").Value = 2687.25
$ws.Range("M34").Value = -1232.25
$ws.Range("N34").Value = -3091.25

$ws.Range("H50").Value = 18647.166
$ws.Range("J50").Value = 18360
$ws.Range("L50").Value = 18360
$ws.Range("N50").Value = -19610

$ws.Range("H60").Value = 12000
$ws.Range("I60").Value = 12000
$ws.Range("K60").Value = 12000
$ws.Range("M60").Value = -11489

$ws.Range("H86").Value = 10023608
$ws.Range("I86").Value = 11454981
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 11454981
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -11453858
$ws.Range("N86").Value = -6246

$ws.Range("H89").Value = 10023608
$ws.Range("I89").Value = 11454981
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 57274905
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -57269289
$ws.Range("N89").Value = -31232

$ws.Range("H132").Value = 2791.6155
$ws.Range("I132").Value = 2791.6155
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8374.8465
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5844.8465
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 1562.5
$ws.Range("I134").Value = 1360.2941
$ws.Range("K134").Value = 4080.8823
$ws.Range("M134").Value = -1545.8823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 277639.5
$ws.Range("I128").Value = 277639.5
$ws.Range("K128").Value = 832918.5
$ws.Range("M128").Value = -827938.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H70").Value = 11000
$ws.Range("I70").Value = 11000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 11000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -10730
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 11000
$ws.Range("I73").Value = 11000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 11000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -10064
$ws.Range("N73").ClearContents()

$ws.Range("H102").Value = 5504
$ws.Range("I102").Value = 3756
$ws.Range("J102").Value = 9000
$ws.Range("K102").Value = 3756
$ws.Range("L102").Value = 9000
$ws.Range("M102").Value = -2134
$ws.Range("N102").Value = -12244

$ws.Range("H132").Value = 3217.8
$ws.Range("I132").Value = 3217.8
$ws.Range("K132").Value = 9653.400000000001
$ws.Range("M132").Value = -7123.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 500
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -1090

$ws.Range("H27").Value = 500
$ws.Range("J27").Value = 500
$ws.Range("L27").Value = 500
$ws.Range("N27").Value = -714

$ws.Range("H46").Value = 4863.636
$ws.Range("I46").Value = 4222.222
$ws.Range("J46").Value = 5307.6924
$ws.Range("K46").Value = 4222.222
$ws.Range("L46").Value = 5307.6924
$ws.Range("M46").Value = -4034.222
$ws.Range("N46").Value = -5683.6924

$ws.Range("H93").Value = 893
$ws.Range("I93").Value = 914.6667
$ws.Range("J93").Value = 849.6667
$ws.Range("K93").Value = 914.6667
$ws.Range("L93").Value = 849.6667
$ws.Range("M93").Value = 333.3333
$ws.Range("N93").Value = -3345.6667

$ws.Range("H100").Value = 2198.2
$ws.Range("I100").Value = 2198.2
$ws.Range("K100").Value = 2198.2
$ws.Range("M100").Value = -1657.2

$ws.Range("H132").Value = 449.75
$ws.Range("I132").Value = 449.75
$ws.Range("K132").Value = 1349.25
$ws.Range("M132").Value = 1180.75

$ws.Range("H136").Value = 3120.7778
$ws.Range("I136").Value = 3120.7778
$ws.Range("K136").Value = 9362.3334
$ws.Range("M136").Value = -6812.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1750.7222
$ws.Range("I132").Value = 1500.8125
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 4502.4375
$ws.Range("L132").Value = 11250
$ws.Range("M132").Value = -1972.4375
$ws.Range("N132").Value = -16310
